$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix style of E17: it currently uses the "applyNumberFormat + center" style (index 4).
# The edit removes that unused style and switches the cell to the plain "center" style (index 1).
$ws.Range("E17").HorizontalAlignment = -4108  # xlCenter

# Add the missing score values for train / test / validation rows
# (these inherit the column's existing "center" style automatically).
$ws.Range("E23").Value = 0.9908
$ws.Range("E24").Value = 0.9857
$ws.Range("E25").Value = 0.9852

# Update the active selection to match the author's final cursor position.
$ws.Range("U11").Select()
